$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsTrans = $wb.Worksheets.Item("Transactions")

# Insert a new column before column N (14) on the Repayment Schedule sheet
$wsRepay.Columns.Item(14).Insert()

# Make Repayment Schedule the active sheet / tab
$wsRepay.Select()
$wsRepay.Range("T4").Select()
